$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: ALC
$ws.Range("H12").Value = 866.3333
$ws.Range("J12").Value = 866.3333
$ws.Range("L12").Value = 866.3333
$ws.Range("N12").Value = -1206.3333

# Row 40: ALC
$ws.Range("H40").Value = 1775.375
$ws.Range("J40").Value = 1400.8
$ws.Range("L40").Value = 1400.8
$ws.Range("N40").Value = -1750.8

# Row 64: ALC
$ws.Range("H64").Value = 4025.1667
$ws.Range("J64").Value = 3930
$ws.Range("L64").Value = 3930
$ws.Range("N64").Value = -4426

# Row 67: ALC
$ws.Range("H67").Value = 4025.1667
$ws.Range("J67").Value = 3930
$ws.Range("L67").Value = 3930
$ws.Range("N67").Value = -5646

# Row 98: ALC
$ws.Range("H98").Value = 8030.3887
$ws.Range("I98").Value = 10352.23
$ws.Range("J98").Value = 1993.6
$ws.Range("K98").Value = 10352.23
$ws.Range("L98").Value = 1993.6
$ws.Range("M98").Value = -8854.23
$ws.Range("N98").Value = -4989.6

# Row 122: ALC
$ws.Range("H122").Value = 8030.3887
$ws.Range("I122").Value = 10352.23
$ws.Range("J122").Value = 1993.6
$ws.Range("K122").Value = 31056.69
$ws.Range("L122").Value = 5980.799999999999
$ws.Range("M122").Value = -28606.69
$ws.Range("N122").Value = -10880.8

# Row 138: ALC
$ws.Range("H138").Value = 2851.6453
$ws.Range("J138").Value = 2841.8691
$ws.Range("L138").Value = 8525.6073
$ws.Range("N138").Value = -18805.6073

$ws = $wb.Worksheets.Item("ARM")
# Row 88: ARM
$ws.Range("H88").Value = 1315.4
$ws.Range("I88").Value = 1463.3334
$ws.Range("J88").Value = 1252
$ws.Range("K88").Value = 1463.3334
$ws.Range("L88").Value = 1252
$ws.Range("M88").Value = -1057.3334
$ws.Range("N88").Value = -2064

# Row 91: ARM
$ws.Range("H91").Value = 1315.4
$ws.Range("I91").Value = 1463.3334
$ws.Range("J91").Value = 1252
$ws.Range("K91").Value = 1463.3334
$ws.Range("L91").Value = 1252
$ws.Range("M91").Value = -59.33339999999998
$ws.Range("N91").Value = -4060

# Row 110: ARM
$ws.Range("H110").Value = 5500
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 5500
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 5500
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -9590

$ws = $wb.Worksheets.Item("BSM")
# Row 86: BSM
$ws.Range("H86").Value = 3210.8572
$ws.Range("I86").Value = 3168.0386
$ws.Range("J86").Value = 3334.5557
$ws.Range("K86").Value = 3168.0386
$ws.Range("L86").Value = 3334.5557
$ws.Range("M86").Value = -2045.0386
$ws.Range("N86").Value = -5580.5557

# Row 89: BSM
$ws.Range("H89").Value = 3210.8572
$ws.Range("I89").Value = 3168.0386
$ws.Range("J89").Value = 3334.5557
$ws.Range("K89").Value = 15840.193
$ws.Range("L89").Value = 16672.7785
$ws.Range("M89").Value = -10224.193
$ws.Range("N89").Value = -27904.7785

# Row 94: BSM
$ws.Range("H94").Value = 22728036
$ws.Range("I94").Value = 25000768
$ws.Range("J94").Value = 699
$ws.Range("K94").Value = 25000768
$ws.Range("L94").Value = 699
$ws.Range("M94").Value = -25000317
$ws.Range("N94").Value = -1601

# Row 105: BSM
$ws.Range("H105").Value = 62501530
$ws.Range("I105").Value = 66668096
$ws.Range("K105").Value = 66668096
$ws.Range("M105").Value = -66666349

# Row 134: BSM
$ws.Range("H134").Value = 3843.389
$ws.Range("I134").Value = 816.5
$ws.Range("J134").Value = 8599.929
$ws.Range("K134").Value = 2449.5
$ws.Range("L134").Value = 25799.787
$ws.Range("M134").Value = 85.5
$ws.Range("N134").Value = -30869.787

$ws = $wb.Worksheets.Item("CRP")
# Row 22: CRP
$ws.Range("H22").Value = 35435
$ws.Range("I22").Value = 548.6667
$ws.Range("J22").Value = 50386.285
$ws.Range("K22").Value = 548.6667
$ws.Range("L22").Value = 50386.285
$ws.Range("M22").Value = -198.6667
$ws.Range("N22").Value = -51086.285

# Row 31: CRP
$ws.Range("H31").Value = 1769.2413
$ws.Range("I31").Value = 1704.625
$ws.Range("J31").Value = 2079.4
$ws.Range("K31").Value = 1704.625
$ws.Range("L31").Value = 2079.4
$ws.Range("M31").Value = -1409.625
$ws.Range("N31").Value = -2669.4

# Row 34: CRP
$ws.Range("H34").Value = 1769.2413
$ws.Range("I34").Value = 1704.625
$ws.Range("J34").Value = 2079.4
$ws.Range("K34").Value = 1704.625
$ws.Range("L34").Value = 2079.4
$ws.Range("M34").Value = -1502.625
$ws.Range("N34").Value = -2483.4

# Row 62: CRP
$ws.Range("H62").Value = 20002410
$ws.Range("I62").Value = 2585.7144
$ws.Range("K62").Value = 2585.7144
$ws.Range("M62").Value = -1961.7144

# Row 65: CRP
$ws.Range("H65").Value = 20002410
$ws.Range("I65").Value = 2585.7144
$ws.Range("K65").Value = 12928.572
$ws.Range("M65").Value = -9808.572

# Row 134: CRP
$ws.Range("H134").Value = 8197962.5
$ws.Range("I134").Value = 1133.5682
$ws.Range("K134").Value = 3400.7046
$ws.Range("M134").Value = -865.7046

$ws = $wb.Worksheets.Item("CUL")
# Row 39: CUL
$ws.Range("H39").Value = 4190.2856
$ws.Range("J39").Value = 4255.3335
$ws.Range("L39").Value = 12766.0005
$ws.Range("N39").Value = -13354.0005

# Row 107: CUL
$ws.Range("H107").Value = 5879.9165
$ws.Range("J107").Value = 9819.571
$ws.Range("L107").Value = 29458.713
$ws.Range("N107").Value = -33298.713

# Row 131: CUL
$ws.Range("H131").Value = 35768876
$ws.Range("I131").Value = 142857570
$ws.Range("K131").Value = 428572710
$ws.Range("M131").Value = -428567670

$ws = $wb.Worksheets.Item("GSM")
# Row 42: GSM
$ws.Range("H42").Value = 41194.4
$ws.Range("J42").Value = 41194.4
$ws.Range("L42").Value = 41194.4
$ws.Range("N42").Value = -42164.4

# Row 47: GSM
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# Row 70: GSM
$ws.Range("I70").Value = 15628744
$ws.Range("J70").Value = 20003610
$ws.Range("K70").Value = 15628744
$ws.Range("L70").Value = 20003610
$ws.Range("M70").Value = -15628474
$ws.Range("N70").Value = -20004150

# Row 73: GSM
$ws.Range("I73").Value = 15628744
$ws.Range("J73").Value = 20003610
$ws.Range("K73").Value = 15628744
$ws.Range("L73").Value = 20003610
$ws.Range("M73").Value = -15627808
$ws.Range("N73").Value = -20005482

# Row 80: GSM
$ws.Range("H80").Value = 6200
$ws.Range("J80").Value = 6200
$ws.Range("L80").Value = 6200
$ws.Range("N80").Value = -8196

# Row 83: GSM
$ws.Range("H83").Value = 6200
$ws.Range("J83").Value = 6200
$ws.Range("L83").Value = 31000
$ws.Range("N83").Value = -40984

# Row 115: GSM
$ws.Range("H115").Value = 41194.4
$ws.Range("J115").Value = 41194.4
$ws.Range("L115").Value = 41194.4
$ws.Range("N115").Value = -43544.4

$ws = $wb.Worksheets.Item("LTW")
# Row 141: LTW
$ws.Range("H141").Value = 48309.617
$ws.Range("J141").Value = 47335.418
$ws.Range("L141").Value = 47335.418
$ws.Range("N141").Value = -57695.418

$ws = $wb.Worksheets.Item("WVR")
# Row 27: WVR
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
